$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Phút hành chính" column (column K) entirely, shifting
# subsequent columns left.
$ws.Columns("K").Delete()

# Update selection to match target state
$ws.Range("M9").Select()
